# Auto-generated edit script: fixes swapped match-data rows in "Lithuania A Lyga"
# by writing the corrected B (match id) and F:AC (match data) values per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 ---
$ws.Range("B27").Value2 = 6227923
$row27 = New-Object "object[,]" 1,24
$row27[0,0] = "Hegelmann Litauen"
$row27[0,1] = "Panevezys"
$row27[0,2] = 2
$row27[0,3] = 4
$row27[0,4] = "A"
$row27[0,5] = 2.6
$row27[0,6] = 3.25
$row27[0,7] = 2.375
$row27[0,8] = 2.625
$row27[0,9] = 4
$row27[0,10] = 2.1
$row27[0,11] = 0.25
$row27[0,12] = 1.85
$row27[0,13] = 1.95
$row27[0,14] = 2.5
$row27[0,15] = 1.775
$row27[0,16] = 2.025
$row27[0,17] = -1
$row27[0,18] = -1
$row27[0,19] = 1.1
$row27[0,20] = -1
$row27[0,21] = 0.95
$row27[0,22] = 0.7749999999999999
$row27[0,23] = -1
$ws.Range("F27:AC27").Value2 = $row27

# --- Row 28 ---
$ws.Range("B28").Value2 = 6227232
$row28 = New-Object "object[,]" 1,24
$row28[0,0] = "FK Siauliai"
$row28[0,1] = "FK Zalgiris Vilnius"
$row28[0,2] = 1
$row28[0,3] = 0
$row28[0,4] = "H"
$row28[0,5] = 4.333
$row28[0,6] = 3.75
$row28[0,7] = 1.615
$row28[0,8] = 4
$row28[0,9] = 3.6
$row28[0,10] = 1.7
$row28[0,11] = 0.75
$row28[0,12] = 1.825
$row28[0,13] = 1.975
$row28[0,14] = 2.25
$row28[0,15] = 1.85
$row28[0,16] = 1.95
$row28[0,17] = 3
$row28[0,18] = -1
$row28[0,19] = -1
$row28[0,20] = 0.825
$row28[0,21] = -1
$row28[0,22] = -1
$row28[0,23] = 0.95
$ws.Range("F28:AC28").Value2 = $row28

# --- Row 125 ---
$ws.Range("B125").Value2 = 6732794
$row125 = New-Object "object[,]" 1,24
$row125[0,0] = "FK Siauliai"
$row125[0,1] = "FK Dziugas Telsiai"
$row125[0,2] = 3
$row125[0,3] = 0
$row125[0,4] = "H"
$row125[0,5] = 1.25
$row125[0,6] = 5
$row125[0,7] = 9
$row125[0,8] = 1.25
$row125[0,9] = 5.25
$row125[0,10] = 9
$row125[0,11] = -1.75
$row125[0,12] = 2
$row125[0,13] = 1.8
$row125[0,14] = 3
$row125[0,15] = 1.975
$row125[0,16] = 1.825
$row125[0,17] = 0.25
$row125[0,18] = -1
$row125[0,19] = -1
$row125[0,20] = 1
$row125[0,21] = -1
$row125[0,22] = 0
$row125[0,23] = -0
$ws.Range("F125:AC125").Value2 = $row125

# --- Row 126 ---
$ws.Range("B126").Value2 = 6732795
$row126 = New-Object "object[,]" 1,24
$row126[0,0] = "Suduva Marijampole"
$row126[0,1] = "Banga Gargzdai"
$row126[0,2] = 1
$row126[0,3] = 0
$row126[0,4] = "H"
$row126[0,5] = 2.15
$row126[0,6] = 3.2
$row126[0,7] = 3
$row126[0,8] = 2.3
$row126[0,9] = 3.2
$row126[0,10] = 2.7
$row126[0,11] = -0.25
$row126[0,12] = 2.05
$row126[0,13] = 1.75
$row126[0,14] = 2.25
$row126[0,15] = 1.9
$row126[0,16] = 1.9
$row126[0,17] = 1.3
$row126[0,18] = -1
$row126[0,19] = -1
$row126[0,20] = 1.05
$row126[0,21] = -1
$row126[0,22] = -1
$row126[0,23] = 0.8999999999999999
$ws.Range("F126:AC126").Value2 = $row126

# --- Row 164 ---
$ws.Range("B164").Value2 = 7326568
$row164 = New-Object "object[,]" 1,24
$row164[0,0] = "Hegelmann Litauen"
$row164[0,1] = "Panevezys"
$row164[0,2] = 0
$row164[0,3] = 0
$row164[0,4] = "D"
$row164[0,5] = 2.375
$row164[0,6] = 3.2
$row164[0,7] = 2.625
$row164[0,8] = 2.7
$row164[0,9] = 3.2
$row164[0,10] = 2.3
$row164[0,11] = 0
$row164[0,12] = 2.05
$row164[0,13] = 1.75
$row164[0,14] = 2.25
$row164[0,15] = 1.875
$row164[0,16] = 1.925
$row164[0,17] = -1
$row164[0,18] = 2.2
$row164[0,19] = -1
$row164[0,20] = 0
$row164[0,21] = -0
$row164[0,22] = -1
$row164[0,23] = 0.925
$ws.Range("F164:AC164").Value2 = $row164

# --- Row 165 ---
$ws.Range("B165").Value2 = 6732827
$row165 = New-Object "object[,]" 1,24
$row165[0,0] = "FK Dziugas Telsiai"
$row165[0,1] = "FK Kauno Zalgiris"
$row165[0,2] = 0
$row165[0,3] = 2
$row165[0,4] = "A"
$row165[0,5] = 6
$row165[0,6] = 3.9
$row165[0,7] = 1.444
$row165[0,8] = 4.75
$row165[0,9] = 3.6
$row165[0,10] = 1.65
$row165[0,11] = 0.75
$row165[0,12] = 1.9
$row165[0,13] = 1.9
$row165[0,14] = 2.5
$row165[0,15] = 1.95
$row165[0,16] = 1.85
$row165[0,17] = -1
$row165[0,18] = -1
$row165[0,19] = 0.6499999999999999
$row165[0,20] = -1
$row165[0,21] = 0.8999999999999999
$row165[0,22] = -1
$row165[0,23] = 0.8500000000000001
$ws.Range("F165:AC165").Value2 = $row165

# --- Row 175 ---
$ws.Range("B175").Value2 = 6732834
$row175 = New-Object "object[,]" 1,24
$row175[0,0] = "Panevezys"
$row175[0,1] = "FK Dziugas Telsiai"
$row175[0,2] = 0
$row175[0,3] = 0
$row175[0,4] = "D"
$row175[0,5] = 1.25
$row175[0,6] = 5.5
$row175[0,7] = 7.5
$row175[0,8] = 1.45
$row175[0,9] = 4.5
$row175[0,10] = 5
$row175[0,11] = -1
$row175[0,12] = 1.775
$row175[0,13] = 2.025
$row175[0,14] = 2.5
$row175[0,15] = 1.875
$row175[0,16] = 1.925
$row175[0,17] = -1
$row175[0,18] = 3.5
$row175[0,19] = -1
$row175[0,20] = -1
$row175[0,21] = 1.025
$row175[0,22] = -1
$row175[0,23] = 0.925
$ws.Range("F175:AC175").Value2 = $row175

# --- Row 176 ---
$ws.Range("B176").Value2 = 6732836
$row176 = New-Object "object[,]" 1,24
$row176[0,0] = "FK Siauliai"
$row176[0,1] = "Banga Gargzdai"
$row176[0,2] = 3
$row176[0,3] = 0
$row176[0,4] = "H"
$row176[0,5] = 1.222
$row176[0,6] = 5.5
$row176[0,7] = 9
$row176[0,8] = 1.363
$row176[0,9] = 4.5
$row176[0,10] = 7
$row176[0,11] = -1.25
$row176[0,12] = 1.9
$row176[0,13] = 1.9
$row176[0,14] = 2.5
$row176[0,15] = 1.975
$row176[0,16] = 1.825
$row176[0,17] = 0.363
$row176[0,18] = -1
$row176[0,19] = -1
$row176[0,20] = 0.8999999999999999
$row176[0,21] = -1
$row176[0,22] = 0.9750000000000001
$row176[0,23] = -1
$ws.Range("F176:AC176").Value2 = $row176

# --- Row 177 ---
$ws.Range("B177").Value2 = 6732837
$row177 = New-Object "object[,]" 1,24
$row177[0,0] = "Suduva Marijampole"
$row177[0,1] = "FK Riteriai"
$row177[0,2] = 0
$row177[0,3] = 3
$row177[0,4] = "A"
$row177[0,5] = 3.6
$row177[0,6] = 3.6
$row177[0,7] = 1.8
$row177[0,8] = 3
$row177[0,9] = 3.6
$row177[0,10] = 2
$row177[0,11] = 0.25
$row177[0,12] = 2
$row177[0,13] = 1.8
$row177[0,14] = 2.5
$row177[0,15] = 1.975
$row177[0,16] = 1.825
$row177[0,17] = -1
$row177[0,18] = -1
$row177[0,19] = 1
$row177[0,20] = -1
$row177[0,21] = 0.8
$row177[0,22] = 0.9750000000000001
$row177[0,23] = -1
$ws.Range("F177:AC177").Value2 = $row177

# --- Row 178 ---
$ws.Range("B178").Value2 = 7465686
$row178 = New-Object "object[,]" 1,24
$row178[0,0] = "FK Kauno Zalgiris"
$row178[0,1] = "Hegelmann Litauen"
$row178[0,2] = 4
$row178[0,3] = 2
$row178[0,4] = "H"
$row178[0,5] = 2.3
$row178[0,6] = 4
$row178[0,7] = 2.3
$row178[0,8] = 2.55
$row178[0,9] = 4
$row178[0,10] = 2.2
$row178[0,11] = 0.25
$row178[0,12] = 1.8
$row178[0,13] = 2
$row178[0,14] = 2.75
$row178[0,15] = 1.85
$row178[0,16] = 1.95
$row178[0,17] = 1.55
$row178[0,18] = -1
$row178[0,19] = -1
$row178[0,20] = 0.8
$row178[0,21] = -1
$row178[0,22] = 0.8500000000000001
$row178[0,23] = -1
$ws.Range("F178:AC178").Value2 = $row178

